$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gem Mine Planner")
$ws.Activate()

# Update the "Set Mine Fee" value (was 150, now 300)
$ws.Range("G21").Value = 300

# Update the active selection to reflect where the editor left off
$ws.Range("G24").Select()
